$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "123Qwe,./7"

$ws.Range("H13").Select()
